$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.803.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "'3.261.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'581.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").Value = "'185.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").Value = "'3.260.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").Value = "'6.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").Value = "'0.413"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "'3.823.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").Value = "'27.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("D16").Value = "'67.788.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "'0.0000170"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Value = "'3.264.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D20").Value = "'13.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'394.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.28%  "

$ws.Range("E22").Value = "  -1.75%  "

$ws.Range("D23").Value = "'71.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").Value = "'0.518"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").Value = "'0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "

$ws.Range("E27").Value = "  -2.58%  "

$ws.Range("D28").Value = "'9.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.60%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").Value = "'5.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.25%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("E33").Value = "  -2.37%  "

$ws.Range("D34").Value = "'1.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.27%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").Value = "'162.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").Value = "'1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.51%  "

$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("D39").Value = "'26.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("D40").Value = "'0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.26%  "

$ws.Range("D41").Value = "'4.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.31%  "

$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("E43").Value = "  -5.09%  "

$ws.Range("D44").Value = "'0.0690"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "'40.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").Value = "'2.615.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("D47").Value = "'24.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.73%  "

$ws.Range("D48").Value = "'334.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").Value = "'0.0279"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "

$ws.Range("E51").Value = "  -0.61%  "
